$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "BA1"
$ws.Range("B8").Value = 120

$ws.Range("C8").Select()
